# Apply updated odds values per the commit diff (Jogos_da_Semana_FlashScore_2025-03-05)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.42
$ws.Range("H3").Value = 3.8
$ws.Range("I3").Value = 8.75
$ws.Range("J3").Value = 1.95
$ws.Range("K3").Value = 2.15
$ws.Range("L3").Value = 7.8
$ws.Range("M3").Value = 1.42
$ws.Range("N3").Value = 2.67
$ws.Range("Q3").Value = 3.9
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.47
$ws.Range("X3").Value = 5.3
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 8.75
$ws.Range("AD3").Value = 8
$ws.Range("AE3").Value = 27
$ws.Range("AF3").Value = 200
$ws.Range("AH3").Value = 16
$ws.Range("AJ3").Value = 28
$ws.Range("AL3").Value = 150
$ws.Range("AM3").Value = 150
$ws.Range("I4").Value = 4.05
$ws.Range("L4").Value = 4.55
$ws.Range("X4").Value = 9.5
$ws.Range("AK4").Value = 75
$ws.Range("I5").Value = 5.2
$ws.Range("K5").Value = 2.15
$ws.Range("M5").Value = 1.33
$ws.Range("N5").Value = 3.05
$ws.Range("O5").Value = 1.98
$ws.Range("P5").Value = 1.75
$ws.Range("S5").Value = 1.39
$ws.Range("T5").Value = 2.77
$ws.Range("V5").Value = 1.83
$ws.Range("W5").Value = 6.1
$ws.Range("X5").Value = 7.5
$ws.Range("Z5").Value = 13
$ws.Range("AE5").Value = 15.5
$ws.Range("AH5").Value = 13
$ws.Range("AI5").Value = 32
$ws.Range("AK5").Value = 110
$ws.Range("O8").Value = 1.5
$ws.Range("U8").Value = 1.41
$ws.Range("V8").Value = 2.62
$ws.Range("G10").Value = 1.57
$ws.Range("H10").Value = 3.6
$ws.Range("I10").Value = 5.6
$ws.Range("J10").Value = 2.18
$ws.Range("K10").Value = 2.12
$ws.Range("L10").Value = 5.3
$ws.Range("P10").Value = 1.85
$ws.Range("Q10").Value = 2.75
$ws.Range("R10").Value = 1.34
$ws.Range("U10").Value = 1.75
$ws.Range("V10").Value = 1.85
$ws.Range("X10").Value = 7.4
$ws.Range("Y10").Value = 7.9
$ws.Range("Z10").Value = 11.75
$ws.Range("AC10").Value = 10.5
$ws.Range("AD10").Value = 7.2
$ws.Range("AE10").Value = 15
$ws.Range("AF10").Value = 65
$ws.Range("G12").Value = 8.5
$ws.Range("H12").Value = 5.1
$ws.Range("J12").Value = 7
$ws.Range("K12").Value = 2.65
$ws.Range("L12").Value = 1.65
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 2.6
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 1.65
$ws.Range("U12").Value = 1.8
$ws.Range("V12").Value = 1.91
$ws.Range("W12").Value = 23
$ws.Range("X12").Value = 50
$ws.Range("AB12").Value = 55
$ws.Range("AC12").Value = 17
$ws.Range("AD12").Value = 9.5
$ws.Range("AE12").Value = 16.5
$ws.Range("AF12").Value = 60
$ws.Range("AG12").Value = 350
$ws.Range("AH12").Value = 7.8
$ws.Range("AI12").Value = 6.3
$ws.Range("AJ12").Value = 7.5
$ws.Range("AK12").Value = 7.1
$ws.Range("AL12").Value = 8.5
$ws.Range("AM12").Value = 18.5
$ws.Range("G13").Value = 1.9
$ws.Range("H13").Value = 3.7
$ws.Range("I13").Value = 3.45
$ws.Range("J13").Value = 2.47
$ws.Range("K13").Value = 2.2
$ws.Range("L13").Value = 3.85
$ws.Range("P13").Value = 1.93
$ws.Range("Q13").Value = 2.62
$ws.Range("U13").Value = 1.65
$ws.Range("V13").Value = 2
$ws.Range("W13").Value = 8.25
$ws.Range("X13").Value = 9.5
$ws.Range("Y13").Value = 8.5
$ws.Range("Z13").Value = 16
$ws.Range("AA13").Value = 14.5
$ws.Range("AC13").Value = 13
$ws.Range("AD13").Value = 7.2
$ws.Range("AE13").Value = 14
$ws.Range("AH13").Value = 11.75
$ws.Range("AI13").Value = 19
$ws.Range("AJ13").Value = 12
$ws.Range("AK13").Value = 45
$ws.Range("AL13").Value = 29
$ws.Range("AM13").Value = 35
$ws.Range("AN14").Value = 1.06
$ws.Range("AO14").Value = 8
$ws.Range("G16").Value = 8.25
$ws.Range("H16").Value = 4.65
$ws.Range("J16").Value = 7.2
$ws.Range("L16").Value = 1.83
$ws.Range("M16").Value = 1.23
$ws.Range("N16").Value = 3.75
$ws.Range("O16").Value = 1.7
$ws.Range("P16").Value = 2.05
$ws.Range("Q16").Value = 2.65
$ws.Range("R16").Value = 1.42
$ws.Range("U16").Value = 2.02
$ws.Range("V16").Value = 1.7
$ws.Range("X16").Value = 55
$ws.Range("Y16").Value = 25
$ws.Range("Z16").Value = 200
$ws.Range("AA16").Value = 100
$ws.Range("AB16").Value = 80
$ws.Range("AE16").Value = 22
$ws.Range("AF16").Value = 110
$ws.Range("AH16").Value = 6.7
$ws.Range("AI16").Value = 6.2
$ws.Range("AJ16").Value = 8.5
$ws.Range("AK16").Value = 8.25
$ws.Range("AL16").Value = 11.25
$ws.Range("AM16").Value = 29
$ws.Range("M17").Value = 1.62
$ws.Range("N17").Value = 2.2
$ws.Range("S18").Value = 1.37
$ws.Range("P19").Value = 1.72
$ws.Range("S19").Value = 1.41
$ws.Range("T19").Value = 2.62
$ws.Range("P20").Value = 1.72
$ws.Range("S20").Value = 1.41
$ws.Range("T20").Value = 2.62
$ws.Range("M21").Value = 1.29
$ws.Range("N21").Value = 3.5
$ws.Range("O21").Value = 1.93
$ws.Range("P21").Value = 1.88
$ws.Range("S21").Value = 1.37
